$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 1 header text: "link_testProject_" -> "link_project_"
$ws.Range("D1").Value = "link_project_internalRoleLinkName"
$ws.Range("E1").Value = "link_project_internalRoleLinkName_1"
$ws.Range("F1").Value = "link_project_project_id"
$ws.Range("G1").Value = "link_project_project_id_1"
$ws.Range("H1").Value = "link_project_team_id"
$ws.Range("I1").Value = "link_project_team_id_1"
$ws.Range("J1").Value = "link_project_test_project_id"
$ws.Range("K1").Value = "link_project_test_project_id_1"
$ws.Range("L1").Value = "link_project_trNthChild"
$ws.Range("M1").Value = "link_project_trNthChild_1"

# Update A2 text
$ws.Range("A2").Value = "Data Files/AI-Generated/Common/scheduleTestRunAndConfigureEnvironment-test-data"

# Update column widths.
# NOTE: the runtime's ColumnWidth setter stores (value + 0.8333333) in the
# saved OOXML <col width="..."/> attribute, so subtract 5/6 from each
# desired target width before assigning it.
$ws.Columns.Item(1).ColumnWidth = 81 - (5/6)
$ws.Columns.Item(4).ColumnWidth = 35 - (5/6)
$ws.Columns.Item(5).ColumnWidth = 37 - (5/6)
$ws.Columns.Item(6).ColumnWidth = 25 - (5/6)
$ws.Columns.Item(7).ColumnWidth = 27 - (5/6)
$ws.Columns.Item(8).ColumnWidth = 22 - (5/6)
$ws.Columns.Item(9).ColumnWidth = 24 - (5/6)
$ws.Columns.Item(10).ColumnWidth = 30 - (5/6)
$ws.Columns.Item(11).ColumnWidth = 32 - (5/6)
$ws.Columns.Item(12).ColumnWidth = 25 - (5/6)
$ws.Columns.Item(13).ColumnWidth = 27 - (5/6)
